# Generate Report for Handback
# Adds a new handback entry (cb642401-ca8e-4ad0-bca4-e0fd58c23cfc) as row 4
# on the "Overview", "zh-cn" and "de-de" worksheets, expanding each sheet's
# table (ListObject) from 3 to 4 data rows and wiring up the corresponding
# hyperlinks.
#
# NOTE: values such as "True" / "False" / "" must be written with a leading
# apostrophe so Excel stores them as text (matching the source data, which
# uses text cells, not native booleans).

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"

$hOverview = $wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb642401ca8e4ad0bca4e0fd58c23cfc00000001/e2e/cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md",
    "",
    "",
    "e2e\cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"
)

$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"

$wsOverview.Range("G4").Value = "2016-08-29 12:47:36"
$wsOverview.Range("G4").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$hZhCnA = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb642401ca8e4ad0bca4e0fd58c23cfc00000001/e2e/cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md",
    "",
    "",
    "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"
)

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"

$wsZhCn.Range("G4").Value = "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.3760a90d2bb712b9cb450f117a3854206dd0586e.zh-cn.xlf"

$wsZhCn.Range("H4").Value = "2016-08-29 12:47:31"
$wsZhCn.Range("H4").NumberFormat = $dateFmt

$hZhCnI = $wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/cb642401ca8e4ad0bca4e0fd58c23cfc00000002/e2e/cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md",
    "",
    "",
    "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"
)

$wsZhCn.Range("J4").Value = "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.3760a90d2bb712b9cb450f117a3854206dd0586e.zh-cn.xlf"

$wsZhCn.Range("K4").Value = "2016-08-29 12:47:48"
$wsZhCn.Range("K4").NumberFormat = $dateFmt

$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$hDeDeA = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cb642401ca8e4ad0bca4e0fd58c23cfc00000001/e2e/cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md",
    "",
    "",
    "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"
)

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"

$wsDeDe.Range("G4").Value = "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.3760a90d2bb712b9cb450f117a3854206dd0586e.de-de.xlf"

$wsDeDe.Range("H4").Value = "2016-08-29 12:47:36"
$wsDeDe.Range("H4").NumberFormat = $dateFmt

$hDeDeI = $wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/cb642401ca8e4ad0bca4e0fd58c23cfc00000003/e2e/cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md",
    "",
    "",
    "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.md"
)

$wsDeDe.Range("J4").Value = "cb642401-ca8e-4ad0-bca4-e0fd58c23cfc.3760a90d2bb712b9cb450f117a3854206dd0586e.de-de.xlf"

$wsDeDe.Range("K4").Value = "2016-08-29 12:47:55"
$wsDeDe.Range("K4").NumberFormat = $dateFmt

$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"
